# Reformat PAF plot code, change exposures for aOR table, and create separate PAF and aOR plots
#
# This script edits the "aOR" worksheet of the workbook:
#  - Merges the "Caveat" column text into the "Exposure" column text
#  - Deletes the now-redundant "Caveat" column from the sheet and its Table
#  - Makes "aOR" the active/selected sheet (instead of "PAF"), with a new selection

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("aOR")
$ws.Activate()

# Update the Exposure column (A) text for the rows that previously carried a
# "Caveat" value, folding that caveat text into the exposure description.
# Written in this order so that the shared-string table ends up populated in
# the same sequence as the target file.
$ws.Range("A2").Value2 = "Smoking during final 2w of pregnancy"
$ws.Range("A5").Value2 = "Not exclusively breast feeding on discharge"
$ws.Range("A8").Value2 = "Prone sleeping position relative to back"
$ws.Range("A4").Value2 = "Prone sleeping position relative to other"

# Remove the "Caveat" column (column C) entirely - the Study/aOR/Lower CI/
# Upper CI columns shift left to take its place.
$tbl = $ws.ListObjects.Item(1)
$ws.Columns.Item(3).Delete()
$tbl.Resize($ws.Range("A1:E10"))

# Restore the table header names (column delete does not rename them itself).
$ws.Range("C1").Value2 = "aOR"
$ws.Range("D1").Value2 = "Lower CI"
$ws.Range("E1").Value2 = "Upper CI"

# Make the aOR sheet the active tab with the new selection, and leave the PAF
# sheet unselected.
$ws.Range("K9").Select()
